# Rename header columns D1:L1 from "evento_N" to "fecha N" on every worksheet.
$wb = $excel.ActiveWorkbook

$newHeaders = @("fecha 1", "fecha 2", "fecha 3", "fecha 4", "fecha 5", "fecha 6", "fecha 7", "fecha 8", "fecha 9")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $col = 4 + $i  # Column D is index 4
        $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
    }
}
